$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.051.19"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.18"
$ws.Range("E3").Value = "  +5.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.99"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5104"
$ws.Range("E7").Value = "  +3.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.78"
$ws.Range("E8").Value = "  +3.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2987"
$ws.Range("E9").Value = "  +7.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06794"
$ws.Range("E10").Value = "  +5.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.909.15"
$ws.Range("E11").Value = "  +5.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.23"
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6968"
$ws.Range("E14").Value = "  +7.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.57"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.868"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.069.20"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008148"
$ws.Range("E18").Value = "  +9.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.97"
$ws.Range("E20").Value = "  +5.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.157.83"
$ws.Range("E21").Value = "  +5.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.812"
$ws.Range("E23").Value = "  +4.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.717"
$ws.Range("E24").Value = "  +6.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.220"
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.86"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.98"
$ws.Range("E27").Value = "  +2.05%  "
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("E29").Value = "  +5.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.393"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.240"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08807"
$ws.Range("E32").Value = "  +5.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.998"
$ws.Range("E33").Value = "  +4.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05064"
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.140"
$ws.Range("E35").Value = "  +4.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7135"
$ws.Range("E36").Value = "  +4.86%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.262"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9641"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01687"
$ws.Range("E41").Value = "  +5.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.162"
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "105.38"
$ws.Range("E43").Value = "  +5.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4296"
$ws.Range("E44").Value = "  +4.76%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.618"
$ws.Range("E46").Value = "  +5.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1275"
$ws.Range("E47").Value = "  +4.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05737"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.16"
$ws.Range("E49").Value = "  +4.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.429"
$ws.Range("E50").Value = "  +3.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3793"
$ws.Range("E51").Value = "  +4.14%  "
